# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets, which carry duplicate data for the first 16 rows.

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value (same updates apply to both sheets)
$updates = @{
    2  = 227
    3  = 262
    4  = 272
    5  = 808
    6  = 260
    7  = 6373
    10 = 108
    11 = 69
    14 = 194
    15 = 491
    16 = 45
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
